$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap columns B ("position") and C ("teamAbbreviation") ---
# Cutting column C and inserting it at column B's position moves its content
# (and its custom width formatting) there in one step, shifting the old
# column B one slot to the right into C.
$ws.Columns.Item(3).Cut()
$ws.Columns.Item(2).Insert()

# --- Updated standings after Saint Tropez ---
# columns: A = points, B = teamAbbreviation, C = position
$data = @(
    @(70, "AUS", 1),
    @(68, "GBR", 2),
    @(67, "NZL", 3),
    @(64, "ESP", 4),
    @(52, "FRA", 5),
    @(45, "CAN", 6),
    @(30, "SUI", 7),
    @(26, "DEN", 8),
    @(22, "ITA", 9),
    @(14, "BRA", 10),
    @(13, "GER", 11),
    @(-8, "USA", 12)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}

# --- Refresh the recorded sort state to match the new two-column range ---
$sf = $ws.Sort.SortFields
$sf.Clear()
$sf.Add2($ws.Range("A2:A13"), $null, 2, $null, $null)
$ws.Sort.SetRange($ws.Range("A2:B13"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# --- Selection change ---
$ws.Range("A2:A3").Select()
$ws.Range("A3").Activate()
